# Scheduled-runner refresh: update market-board derived profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns) across several
# sheets of the Ixion_Profits workbook. Values below are the refreshed
# numbers pulled by the runner for the affected leve rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Cells.Item(70, 8).Value = 2989.182
$ws.Cells.Item(70, 9).Value = 5751
$ws.Cells.Item(70, 10).Value = 2375.4443
$ws.Cells.Item(70, 11).Value = 17253
$ws.Cells.Item(70, 12).Value = 7126.3329
$ws.Cells.Item(70, 13).Value = -16983
$ws.Cells.Item(70, 14).Value = -7666.3329
# Row 73
$ws.Cells.Item(73, 8).Value = 2989.182
$ws.Cells.Item(73, 9).Value = 5751
$ws.Cells.Item(73, 10).Value = 2375.4443
$ws.Cells.Item(73, 11).Value = 17253
$ws.Cells.Item(73, 12).Value = 7126.3329
$ws.Cells.Item(73, 13).Value = -16317
$ws.Cells.Item(73, 14).Value = -8998.332900000001
# Row 80
$ws.Cells.Item(80, 8).Value = 533.1111
$ws.Cells.Item(80, 9).Value = 482.8889
$ws.Cells.Item(80, 10).Value = 583.3333
$ws.Cells.Item(80, 11).Value = 1448.6667
$ws.Cells.Item(80, 12).Value = 1749.9999
$ws.Cells.Item(80, 13).Value = -450.6667
$ws.Cells.Item(80, 14).Value = -3745.9999
# Row 82
$ws.Cells.Item(82, 8).Value = 2491.5
$ws.Cells.Item(82, 9).Value = 2491.5
$ws.Cells.Item(82, 11).Value = 7474.5
$ws.Cells.Item(82, 13).Value = -7068.5
# Row 83
$ws.Cells.Item(83, 8).Value = 533.1111
$ws.Cells.Item(83, 9).Value = 482.8889
$ws.Cells.Item(83, 10).Value = 583.3333
$ws.Cells.Item(83, 11).Value = 4346.0001
$ws.Cells.Item(83, 12).Value = 5249.9997
$ws.Cells.Item(83, 13).Value = 645.9998999999998
$ws.Cells.Item(83, 14).Value = -15233.9997
# Row 85
$ws.Cells.Item(85, 8).Value = 2491.5
$ws.Cells.Item(85, 9).Value = 2491.5
$ws.Cells.Item(85, 11).Value = 7474.5
$ws.Cells.Item(85, 13).Value = -6070.5
# Row 129
$ws.Cells.Item(129, 8).Value = 1062.1343
$ws.Cells.Item(129, 9).Value = 773.1
$ws.Cells.Item(129, 10).Value = 1112.8422
$ws.Cells.Item(129, 11).Value = 2319.3
$ws.Cells.Item(129, 12).Value = 3338.5266
$ws.Cells.Item(129, 13).Value = 2680.7
$ws.Cells.Item(129, 14).Value = -13338.5266
# Row 132
$ws.Cells.Item(132, 8).Value = 962.1539
$ws.Cells.Item(132, 9).Value = 790.7347
$ws.Cells.Item(132, 10).Value = 3762
$ws.Cells.Item(132, 11).Value = 2372.2041
$ws.Cells.Item(132, 12).Value = 11286
$ws.Cells.Item(132, 13).Value = 157.7959000000001
$ws.Cells.Item(132, 14).Value = -16346
# Row 138
$ws.Cells.Item(138, 8).Value = 973.15
$ws.Cells.Item(138, 9).Value = 468.36
$ws.Cells.Item(138, 10).Value = 1477.94
$ws.Cells.Item(138, 11).Value = 1405.08
$ws.Cells.Item(138, 12).Value = 4433.82
$ws.Cells.Item(138, 13).Value = 3734.92
$ws.Cells.Item(138, 14).Value = -14713.82
# Row 141
$ws.Cells.Item(141, 8).Value = 1397.6182
$ws.Cells.Item(141, 9).Value = 959.15
$ws.Cells.Item(141, 11).Value = 2877.45
$ws.Cells.Item(141, 13).Value = 2302.55

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 6852.653
$ws.Cells.Item(32, 9).Value = 5689.74
$ws.Cells.Item(32, 10).Value = 11116.667
$ws.Cells.Item(32, 11).Value = 5689.74
$ws.Cells.Item(32, 12).Value = 11116.667
$ws.Cells.Item(32, 13).Value = -5402.74
$ws.Cells.Item(32, 14).Value = -11690.667
# Row 74
$ws.Cells.Item(74, 8).Value = 1176.1786
$ws.Cells.Item(74, 9).Value = 942.75
$ws.Cells.Item(74, 10).Value = 2576.75
$ws.Cells.Item(74, 11).Value = 942.75
$ws.Cells.Item(74, 12).Value = 2576.75
$ws.Cells.Item(74, 13).Value = -68.75
$ws.Cells.Item(74, 14).Value = -4324.75
# Row 77
$ws.Cells.Item(77, 8).Value = 1176.1786
$ws.Cells.Item(77, 9).Value = 942.75
$ws.Cells.Item(77, 10).Value = 2576.75
$ws.Cells.Item(77, 11).Value = 4713.75
$ws.Cells.Item(77, 12).Value = 12883.75
$ws.Cells.Item(77, 13).Value = -345.75
$ws.Cells.Item(77, 14).Value = -21619.75
# Row 132
$ws.Cells.Item(132, 8).Value = 2389.9194
$ws.Cells.Item(132, 9).Value = 1627.4894
$ws.Cells.Item(132, 10).Value = 4778.8667
$ws.Cells.Item(132, 11).Value = 4882.468199999999
$ws.Cells.Item(132, 12).Value = 14336.6001
$ws.Cells.Item(132, 13).Value = -2352.468199999999
$ws.Cells.Item(132, 14).Value = -19396.6001

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Cells.Item(58, 8).Value = 1345.2037
$ws.Cells.Item(58, 9).Value = 747.28125
$ws.Cells.Item(58, 10).Value = 2214.9092
$ws.Cells.Item(58, 11).Value = 747.28125
$ws.Cells.Item(58, 12).Value = 2214.9092
$ws.Cells.Item(58, 13).Value = -544.28125
$ws.Cells.Item(58, 14).Value = -2620.9092
# Row 86
$ws.Cells.Item(86, 8).Value = 2018.3334
$ws.Cells.Item(86, 9).Value = 2023.7778
$ws.Cells.Item(86, 10).Value = 2012.8889
$ws.Cells.Item(86, 11).Value = 2023.7778
$ws.Cells.Item(86, 12).Value = 2012.8889
$ws.Cells.Item(86, 13).Value = -900.7778000000001
$ws.Cells.Item(86, 14).Value = -4258.8889
# Row 89
$ws.Cells.Item(89, 8).Value = 2018.3334
$ws.Cells.Item(89, 9).Value = 2023.7778
$ws.Cells.Item(89, 10).Value = 2012.8889
$ws.Cells.Item(89, 11).Value = 10118.889
$ws.Cells.Item(89, 12).Value = 10064.4445
$ws.Cells.Item(89, 13).Value = -4502.889000000001
$ws.Cells.Item(89, 14).Value = -21296.4445
# Row 132
$ws.Cells.Item(132, 8).Value = 1795.7361
$ws.Cells.Item(132, 9).Value = 1354.3036
$ws.Cells.Item(132, 10).Value = 3340.75
$ws.Cells.Item(132, 11).Value = 4062.9108
$ws.Cells.Item(132, 12).Value = 10022.25
$ws.Cells.Item(132, 13).Value = -1532.9108
$ws.Cells.Item(132, 14).Value = -15082.25
# Row 136
$ws.Cells.Item(136, 8).Value = 1345.2037
$ws.Cells.Item(136, 9).Value = 747.28125
$ws.Cells.Item(136, 10).Value = 2214.9092
$ws.Cells.Item(136, 11).Value = 2241.84375
$ws.Cells.Item(136, 12).Value = 6644.7276
$ws.Cells.Item(136, 13).Value = 308.15625
$ws.Cells.Item(136, 14).Value = -11744.7276

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 286797.6
$ws.Cells.Item(5, 9).Value = 422
$ws.Cells.Item(5, 10).Value = 456019.53
$ws.Cells.Item(5, 11).Value = 1266
$ws.Cells.Item(5, 12).Value = 1368058.59
$ws.Cells.Item(5, 13).Value = -1154
$ws.Cells.Item(5, 14).Value = -1368282.59
# Row 113
$ws.Cells.Item(113, 8).Value = 2069563.1
$ws.Cells.Item(113, 9).Value = 3333949.2
$ws.Cells.Item(113, 10).Value = 714863.5600000001
$ws.Cells.Item(113, 11).Value = 10001847.6
$ws.Cells.Item(113, 12).Value = 2144590.68
$ws.Cells.Item(113, 13).Value = -9999677.600000001
$ws.Cells.Item(113, 14).Value = -2148930.68
# Row 127
$ws.Cells.Item(127, 8).Value = 5682451
$ws.Cells.Item(127, 10).Value = 5682451
$ws.Cells.Item(127, 12).Value = 17047353
$ws.Cells.Item(127, 14).Value = -17057273
# Row 131
$ws.Cells.Item(131, 8).Value = 2223012.8
$ws.Cells.Item(131, 10).Value = 1021.5484
$ws.Cells.Item(131, 12).Value = 3064.6452
$ws.Cells.Item(131, 14).Value = -13144.6452
# Row 135
$ws.Cells.Item(135, 8).Value = 286797.6
$ws.Cells.Item(135, 9).Value = 422
$ws.Cells.Item(135, 10).Value = 456019.53
$ws.Cells.Item(135, 11).Value = 3798
$ws.Cells.Item(135, 12).Value = 4104175.77
$ws.Cells.Item(135, 13).Value = -1263
$ws.Cells.Item(135, 14).Value = -4109245.77

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Cells.Item(132, 8).Value = 1981.8727
$ws.Cells.Item(132, 9).Value = 1558.1282
$ws.Cells.Item(132, 10).Value = 3014.75
$ws.Cells.Item(132, 11).Value = 4674.3846
$ws.Cells.Item(132, 12).Value = 9044.25
$ws.Cells.Item(132, 13).Value = -2144.3846
$ws.Cells.Item(132, 14).Value = -14104.25

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 1216.4615
$ws.Cells.Item(81, 9).Value = 1114.4445
$ws.Cells.Item(81, 10).Value = 1446
$ws.Cells.Item(81, 11).Value = 2228.889
$ws.Cells.Item(81, 12).Value = 2892
$ws.Cells.Item(81, 13).Value = -1167.889
$ws.Cells.Item(81, 14).Value = -5014
# Row 84
$ws.Cells.Item(84, 8).Value = 1216.4615
$ws.Cells.Item(84, 9).Value = 1114.4445
$ws.Cells.Item(84, 10).Value = 1446
$ws.Cells.Item(84, 11).Value = 11144.445
$ws.Cells.Item(84, 12).Value = 14460
$ws.Cells.Item(84, 13).Value = -5840.445
$ws.Cells.Item(84, 14).Value = -25068
# Row 101
$ws.Cells.Item(101, 8).Value = 19801
$ws.Cells.Item(101, 10).Value = 19801
$ws.Cells.Item(101, 12).Value = 19801
$ws.Cells.Item(101, 14).Value = -26291
# Row 132
$ws.Cells.Item(132, 8).Value = 16787.3
$ws.Cells.Item(132, 9).Value = 18521.732
$ws.Cells.Item(132, 10).Value = 2911.8572
$ws.Cells.Item(132, 11).Value = 55565.196
$ws.Cells.Item(132, 12).Value = 8735.571599999999
$ws.Cells.Item(132, 13).Value = -53035.196
$ws.Cells.Item(132, 14).Value = -13795.5716
# Row 136
$ws.Cells.Item(136, 8).Value = 9806931
$ws.Cells.Item(136, 9).Value = 3706.9285
$ws.Cells.Item(136, 10).Value = 21741292
$ws.Cells.Item(136, 11).Value = 11120.7855
$ws.Cells.Item(136, 12).Value = 65223876
$ws.Cells.Item(136, 13).Value = -8570.7855
$ws.Cells.Item(136, 14).Value = -65228976

